# Apply scheduled-runner profit/price updates to the Chocobo_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6186
$ws.Range("I74").Value = 4700.4
$ws.Range("K74").Value = 4700.4
$ws.Range("M74").Value = -3764.4
$ws.Range("H77").Value = 6186
$ws.Range("I77").Value = 4700.4
$ws.Range("K77").Value = 23502
$ws.Range("M77").Value = -18822
$ws.Range("H113").Value = 4124.125
$ws.Range("J113").Value = 4730.846
$ws.Range("L113").Value = 4730.846
$ws.Range("N113").Value = -11238.846
$ws.Range("H116").Value = 391573.78
$ws.Range("J116").Value = 9854.200000000001
$ws.Range("L116").Value = 9854.200000000001
$ws.Range("N116").Value = -16738.2
$ws.Range("H132").Value = 37409076
$ws.Range("I132").Value = 40001204
$ws.Range("J132").Value = 5007503
$ws.Range("K132").Value = 120003612
$ws.Range("L132").Value = 15022509
$ws.Range("M132").Value = -120001082
$ws.Range("N132").Value = -15027569
$ws.Range("H137").Value = 1163326.2
$ws.Range("I137").Value = 2071793.4
$ws.Range("J137").Value = 2507.0557
$ws.Range("K137").Value = 6215380.199999999
$ws.Range("L137").Value = 7521.1671
$ws.Range("M137").Value = -6212830.199999999
$ws.Range("N137").Value = -12621.1671
$ws.Range("H138").Value = 4748.2
$ws.Range("I138").Value = 838.28
$ws.Range("J138").Value = 6051.507
$ws.Range("K138").Value = 2514.84
$ws.Range("L138").Value = 18154.521
$ws.Range("M138").Value = 2625.16
$ws.Range("N138").Value = -28434.521

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3302.2444
$ws.Range("I61").Value = 1020.05554
$ws.Range("K61").Value = 1020.05554
$ws.Range("M61").Value = -808.05554
$ws.Range("H74").Value = 3902.3142
$ws.Range("I74").Value = 4631.905
$ws.Range("K74").Value = 4631.905
$ws.Range("M74").Value = -3757.905
$ws.Range("H77").Value = 3902.3142
$ws.Range("I77").Value = 4631.905
$ws.Range("K77").Value = 23159.525
$ws.Range("M77").Value = -18791.525
$ws.Range("H132").Value = 1367.7662
$ws.Range("I132").Value = 853.1799999999999
$ws.Range("K132").Value = 2559.54
$ws.Range("M132").Value = -29.53999999999996
$ws.Range("H136").Value = 3302.2444
$ws.Range("I136").Value = 1020.05554
$ws.Range("K136").Value = 3060.16662
$ws.Range("M136").Value = -510.16662
$ws.Range("H137").Value = 45195
$ws.Range("J137").Value = 45195
$ws.Range("L137").Value = 45195
$ws.Range("N137").Value = -55395

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1916.6666
$ws.Range("I86").Value = 1916.6666
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1916.6666
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -793.6666
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1916.6666
$ws.Range("I89").Value = 1916.6666
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9583.333000000001
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -3967.333000000001
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 3616.5881
$ws.Range("I99").Value = 1431.8334
$ws.Range("J99").Value = 4808.273
$ws.Range("K99").Value = 1431.8334
$ws.Range("L99").Value = 4808.273
$ws.Range("M99").Value = 66.16660000000002
$ws.Range("N99").Value = -7804.273
$ws.Range("H134").Value = 4300.0684
$ws.Range("J134").Value = 11161.77
$ws.Range("L134").Value = 33485.31
$ws.Range("N134").Value = -38555.31

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 411.16666
$ws.Range("I7").Value = 371.22223
$ws.Range("K7").Value = 371.22223
$ws.Range("M7").Value = -258.22223
$ws.Range("H31").Value = 2191.762
$ws.Range("I31").Value = 955.6667
$ws.Range("J31").Value = 3839.889
$ws.Range("K31").Value = 955.6667
$ws.Range("L31").Value = 3839.889
$ws.Range("M31").Value = -660.6667
$ws.Range("N31").Value = -4429.889
$ws.Range("H34").Value = 2191.762
$ws.Range("I34").Value = 955.6667
$ws.Range("J34").Value = 3839.889
$ws.Range("K34").Value = 955.6667
$ws.Range("L34").Value = 3839.889
$ws.Range("M34").Value = -753.6667
$ws.Range("N34").Value = -4243.889
$ws.Range("H58").Value = 2460.8125
$ws.Range("I58").Value = 1520.3881
$ws.Range("J58").Value = 7307.615
$ws.Range("K58").Value = 1520.3881
$ws.Range("L58").Value = 7307.615
$ws.Range("M58").Value = -1317.3881
$ws.Range("N58").Value = -7713.615
$ws.Range("H99").Value = 3569.9048
$ws.Range("I99").Value = 1715.7
$ws.Range("J99").Value = 5255.5454
$ws.Range("K99").Value = 1715.7
$ws.Range("L99").Value = 5255.5454
$ws.Range("M99").Value = -217.7
$ws.Range("N99").Value = -8251.545399999999
$ws.Range("H126").Value = 3569.9048
$ws.Range("I126").Value = 1715.7
$ws.Range("J126").Value = 5255.5454
$ws.Range("K126").Value = 5147.1
$ws.Range("L126").Value = 15766.6362
$ws.Range("M126").Value = -2677.1
$ws.Range("N126").Value = -20706.6362
$ws.Range("H132").Value = 2643.3953
$ws.Range("I132").Value = 2201.7896
$ws.Range("K132").Value = 6605.3688
$ws.Range("M132").Value = -4075.3688
$ws.Range("H134").Value = 2091.8635
$ws.Range("I134").Value = 1132.5625
$ws.Range("K134").Value = 3397.6875
$ws.Range("M134").Value = -862.6875
$ws.Range("H136").Value = 2460.8125
$ws.Range("I136").Value = 1520.3881
$ws.Range("J136").Value = 7307.615
$ws.Range("K136").Value = 4561.164299999999
$ws.Range("L136").Value = 21922.845
$ws.Range("M136").Value = -2011.164299999999
$ws.Range("N136").Value = -27022.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 9100
$ws.Range("I74").Value = 8500
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 25500
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = -24439
$ws.Range("N74").Value = -32122
$ws.Range("H77").Value = 9100
$ws.Range("I77").Value = 8500
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 76500
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = -71196
$ws.Range("N77").Value = -100608
$ws.Range("H107").Value = 72900
$ws.Range("I107").Value = 539.9
$ws.Range("K107").Value = 1619.7
$ws.Range("M107").Value = 300.3000000000002
$ws.Range("H113").Value = 577.6667
$ws.Range("I113").Value = 583.3
$ws.Range("J113").Value = 566.4
$ws.Range("K113").Value = 1749.9
$ws.Range("L113").Value = 1699.2
$ws.Range("M113").Value = 420.1000000000001
$ws.Range("N113").Value = -6039.2
$ws.Range("H129").Value = 2405.25
$ws.Range("I129").Value = 2278.3076
$ws.Range("K129").Value = 6834.9228
$ws.Range("M129").Value = -1834.9228
$ws.Range("H131").Value = 771.45
$ws.Range("I131").Value = 375.7143
$ws.Range("J131").Value = 801.2366
$ws.Range("K131").Value = 1127.1429
$ws.Range("L131").Value = 2403.7098
$ws.Range("M131").Value = 3912.8571
$ws.Range("N131").Value = -12483.7098
$ws.Range("H138").Value = 3198
$ws.Range("J138").Value = 3216.6667
$ws.Range("L138").Value = 9650.000100000001
$ws.Range("N138").Value = -19930.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2632.8262
$ws.Range("I132").Value = 1464
$ws.Range("J132").Value = 3531.923
$ws.Range("K132").Value = 4392
$ws.Range("L132").Value = 10595.769
$ws.Range("M132").Value = -1862
$ws.Range("N132").Value = -15655.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5064.3
$ws.Range("I40").Value = 4366.5
$ws.Range("J40").Value = 9600
$ws.Range("K40").Value = 4366.5
$ws.Range("L40").Value = 9600
$ws.Range("M40").Value = -4230.5
$ws.Range("N40").Value = -9872
$ws.Range("H132").Value = 13570.233
$ws.Range("I132").Value = 17537.6
$ws.Range("K132").Value = 52612.8
$ws.Range("M132").Value = -50082.8
$ws.Range("H136").Value = 3200.3784
$ws.Range("I136").Value = 1559.7142
$ws.Range("J136").Value = 5353.75
$ws.Range("K136").Value = 4679.142599999999
$ws.Range("L136").Value = 16061.25
$ws.Range("M136").Value = -2129.142599999999
$ws.Range("N136").Value = -21161.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1661.1428
$ws.Range("I132").Value = 1105.2941
$ws.Range("J132").Value = 4023.5
$ws.Range("K132").Value = 3315.8823
$ws.Range("L132").Value = 12070.5
$ws.Range("M132").Value = -785.8823000000002
$ws.Range("N132").Value = -17130.5
$ws.Range("H136").Value = 2722.762
$ws.Range("I136").Value = 1830.7742
$ws.Range("K136").Value = 5492.3226
$ws.Range("M136").Value = -2942.3226
